# Applies the cryptos-list refresh described in the commit diff:
# updates Price (D) and Volume(1h) (E) text cells for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text numbers (e.g. "76.336.53"); force text format
# while writing so Excel does not auto-convert them to numeric values,
# then restore the default "Normal" style so no stray formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "76.336.53"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "2.872.90"
$ws.Range("E3").Value = "  +7.88%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "195.62"
$ws.Range("E5").Value = "  +4.73%  "
$ws.Range("D6").Value = "598.05"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +4.03%  "
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").Value = "2.870.68"
$ws.Range("E10").Value = "  +7.78%  "
$ws.Range("E11").Value = "  +10.31%  "
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("E13").Value = "  +4.27%  "
$ws.Range("D14").Value = "3.398.17"
$ws.Range("E14").Value = "  +7.66%  "
$ws.Range("D15").Value = "76.174.39"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "27.49"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "2.874.96"
$ws.Range("E18").Value = "  +7.31%  "
$ws.Range("D20").Value = "12.52"
$ws.Range("E20").Value = "  +5.05%  "
$ws.Range("D21").Value = "381.65"
$ws.Range("E21").Value = "  +3.01%  "
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("D24").Value = "71.58"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "3.028.40"
$ws.Range("E26").Value = "  +8.35%  "
$ws.Range("D27").Value = "4.22"
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("E28").Value = "  +4.87%  "
$ws.Range("D29").Value = "0.0000105"
$ws.Range("E29").Value = "  +11.01%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "509.50"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "7.72"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "167.20"
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("E37").Value = "  +4.58%  "
$ws.Range("D38").Value = "0.117"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "19.52"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").Value = "184.61"
$ws.Range("E40").Value = "  +9.22%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "0.345"
$ws.Range("E42").Value = "  +5.02%  "
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("D44").Value = "1.68"
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("E45").Value = "  +9.25%  "
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("D47").Value = "40.27"
$ws.Range("E47").Value = "  +2.99%  "
$ws.Range("D48").Value = "2.36"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "0.685"
$ws.Range("E49").Value = "  +16.43%  "
$ws.Range("D50").Value = "0.579"
$ws.Range("E50").Value = "  +8.84%  "
$ws.Range("D51").Value = "3.75"
$ws.Range("E51").Value = "  +3.01%  "

$ws.Range("D2:D51").Style = "Normal"

